$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3, shifting old row 3 down to row 4
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value = "003_Profil_002_Profil_Datenaenderung"
$ws.Range("B3").Value = "var003_Profil_002_Profil_Datenaenderung"
$ws.Range("C3").Value = "001_Login_001_Successful"
$ws.Range("D3").Value = "Menueauswahl Mein Profil"
$ws.Range("E3").Value = "003_Profil_002_Profil_Datenaenderung"
$ws.Range("F3").Value = "Logoff"

$ws.Range("E10").Select()
